$wb = $excel.ActiveWorkbook

# --- "axes" sheet: insert 3 new columns (D:F) for arrow labels + w/w units,
#     pushing the previous Title/Folk-classification column out to G.
$ws = $wb.Worksheets.Item("axes")
$ws.Columns("D:F").Insert()

$ws.Range("D1").Value = "A_arrow"
$ws.Range("E1").Value = "B_arrow"
$ws.Range("F1").Value = "C_arrow"

# Populate row 2 in shared-string creation order Silt, Clay, Sand so the
# workbook's shared-strings table matches the authored order.
$ws.Range("F2").Value = "Silt (w/w)"
$ws.Range("E2").Value = "Clay (w/w)"
$ws.Range("D2").Value = "Sand (w/w)"

$ws.Columns("D:F").ColumnWidth = 14.08

# This sheet becomes the active tab/selection on reopen.
$ws.Activate()
[void]$ws.Range("E7").Select()
